$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 369, shifting existing rows 369:379 down to 372:382
$ws.Rows("369:371").Insert()

# New row 369: Mandarina, Clemenuless, Especial
$ws.Range("A369").Value = 2
$ws.Range("B369").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C369").Value = "Coquimbo"
$ws.Range("D369").Value = 44714
$ws.Range("D369").NumberFormat = $ws.Range("D372").NumberFormat
$ws.Range("E369").Value = 4
$ws.Range("F369").Value = "Fruta"
$ws.Range("G369").Value = 100102
$ws.Range("H369").Value = "Cítricos"
$ws.Range("I369").Value = 100102004
$ws.Range("J369").Value = "Mandarina"
$ws.Range("K369").Value = "Clemenuless"
$ws.Range("L369").Value = "Especial"
$ws.Range("M369").Value = 20
$ws.Range("N369").Value = 235000
$ws.Range("O369").Value = 240000
$ws.Range("P369").Value = 237500
$ws.Range("Q369").Value = "`$/bins (450 kilos)"
$ws.Range("R369").Value = "Provincia de Limarí"
$ws.Range("S369").Value = 528
$ws.Range("T369").Value = 450

# New row 370: Mandarina, Clemenuless, Primera
$ws.Range("A370").Value = 2
$ws.Range("B370").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C370").Value = "Coquimbo"
$ws.Range("D370").Value = 44714
$ws.Range("D370").NumberFormat = $ws.Range("D372").NumberFormat
$ws.Range("E370").Value = 4
$ws.Range("F370").Value = "Fruta"
$ws.Range("G370").Value = 100102
$ws.Range("H370").Value = "Cítricos"
$ws.Range("I370").Value = 100102004
$ws.Range("J370").Value = "Mandarina"
$ws.Range("K370").Value = "Clemenuless"
$ws.Range("L370").Value = "Primera"
$ws.Range("M370").Value = 20
$ws.Range("N370").Value = 195000
$ws.Range("O370").Value = 200000
$ws.Range("P370").Value = 197500
$ws.Range("Q370").Value = "`$/bins (450 kilos)"
$ws.Range("R370").Value = "Provincia de Limarí"
$ws.Range("S370").Value = 439
$ws.Range("T370").Value = 450

# New row 371: Mandarina, Clemenuless, Segunda
$ws.Range("A371").Value = 2
$ws.Range("B371").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C371").Value = "Coquimbo"
$ws.Range("D371").Value = 44714
$ws.Range("D371").NumberFormat = $ws.Range("D372").NumberFormat
$ws.Range("E371").Value = 4
$ws.Range("F371").Value = "Fruta"
$ws.Range("G371").Value = 100102
$ws.Range("H371").Value = "Cítricos"
$ws.Range("I371").Value = 100102004
$ws.Range("J371").Value = "Mandarina"
$ws.Range("K371").Value = "Clemenuless"
$ws.Range("L371").Value = "Segunda"
$ws.Range("M371").Value = 16
$ws.Range("N371").Value = 155000
$ws.Range("O371").Value = 160000
$ws.Range("P371").Value = 157500
$ws.Range("Q371").Value = "`$/bins (450 kilos)"
$ws.Range("R371").Value = "Provincia de Limarí"
$ws.Range("S371").Value = 350
$ws.Range("T371").Value = 450
